# Updates cryptocurrency price/volume data in the worksheet to reflect the
# latest snapshot (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.807.05"
$ws.Range("E2").Value = '  -0.32%  '

$ws.Range("D3").Value = "'1.636.03"
$ws.Range("E3").Value = '  +0.01%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = "'215.68"
$ws.Range("E5").Value = '  +0.62%  '

$ws.Range("D6").Value = "'0.5056"
$ws.Range("E6").Value = '  -0.20%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").Value = "'0.2574"
$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").Value = "'0.06422"
$ws.Range("E9").Value = '  +1.02%  '

$ws.Range("D10").Value = "'19.55"
$ws.Range("E10").Value = '  -0.87%  '

$ws.Range("D11").Value = "'0.07780"
$ws.Range("E11").Value = '  +0.35%  '

$ws.Range("E12").Value = '  -0.32%  '

$ws.Range("D13").Value = "'1.861.97"
$ws.Range("E13").Value = '  +0.01%  '

$ws.Range("D14").Value = "'1.635.09"
$ws.Range("E14").Value = '  +0.28%  '

$ws.Range("D15").Value = "'0.5628"
$ws.Range("E15").Value = '  +3.29%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = "'0.0₅7590"
$ws.Range("E16").Value = '  -1.87%  '

$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = "'63.12"
$ws.Range("E17").Value = '  -1.53%  '

$ws.Range("D18").Value = "'25.846.51"
$ws.Range("E18").Value = '  -0.27%  '

$ws.Range("D19").Value = "'1.003"
$ws.Range("E19").Value = '  +0.14%  '

$ws.Range("D20").Value = "'195.03"
$ws.Range("E20").Value = '  -0.32%  '

$ws.Range("D21").Value = "'4.319"
$ws.Range("E21").Value = '  -2.73%  '

$ws.Range("D22").Value = "'9.870"
$ws.Range("E22").Value = '  -0.49%  '

$ws.Range("D23").Value = "'6.099"
$ws.Range("E23").Value = '  -0.52%  '

$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("D25").Value = "'1.798"
$ws.Range("E25").Value = '  -4.75%  '

$ws.Range("D26").Value = "'0.1273"
$ws.Range("E26").Value = '  +2.27%  '

$ws.Range("D27").Value = "'139.87"
$ws.Range("E27").Value = '  -2.15%  '

$ws.Range("D28").Value = "'6.773"
$ws.Range("E28").Value = '  -0.90%  '

$ws.Range("E29").Value = '  -1.07%  '

$ws.Range("D30").Value = "'1.242"
$ws.Range("E30").Value = '  +0.43%  '

$ws.Range("D31").Value = "'0.04867"
$ws.Range("E31").Value = '  -0.11%  '

$ws.Range("D32").Value = "'3.293"
$ws.Range("E32").Value = '  +1.63%  '

$ws.Range("D33").Value = "'3.214"
$ws.Range("E33").Value = '  +0.56%  '

$ws.Range("D34").Value = "'1.555"
$ws.Range("E34").Value = '  +0.42%  '

$ws.Range("D35").Value = "'2.374"
$ws.Range("E35").Value = '  +0.15%  '

$ws.Range("D36").Value = "'0.9026"
$ws.Range("E36").Value = '  -1.07%  '

$ws.Range("D37").Value = "'2.574"
$ws.Range("E37").Value = '  +0.17%  '

$ws.Range("D38").Value = "'1.131.07"
$ws.Range("E38").Value = '  +0.80%  '

$ws.Range("D39").Value = "'0.5500"
$ws.Range("E39").Value = '  -0.29%  '

$ws.Range("D40").Value = "'0.01561"
$ws.Range("E40").Value = '  -0.20%  '

$ws.Range("D41").Value = "'0.9953"

$ws.Range("D42").Value = "'5.519"
$ws.Range("E42").Value = '  -1.24%  '

$ws.Range("D43").Value = "'0.8006"
$ws.Range("E43").Value = '  -0.45%  '

$ws.Range("D44").Value = "'97.77"
$ws.Range("E44").Value = '  -0.79%  '

$ws.Range("D45").Value = "'1.772.62"
$ws.Range("E45").Value = '  +0.12%  '

$ws.Range("E46").Value = '  -7.77%  '

$ws.Range("D47").Value = "'55.33"
$ws.Range("E47").Value = '  +0.65%  '

$ws.Range("D48").Value = "'0.4390"
$ws.Range("E48").Value = '  -2.04%  '

$ws.Range("D49").Value = "'0.05052"
$ws.Range("E49").Value = '  -2.50%  '

$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = '  -0.25%  '

